# Restructure "Hoja1" from a 5-column CODIGO/NOMBRE/CAPACITACION/FECHA/URL
# sheet (with two rows of data plus a stray A3 cell) into a 3-column
# NOMBRE/CAPACITACION/URL sheet with four identical data rows (2-5), no
# hyperlink objects (the URL text itself stays, but is no longer a live
# hyperlink), and new column widths/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two hyperlink objects that used to live on A2 and E2.
$ws.Hyperlinks.Delete()

# The old D and E columns (FECHA / URL-link) go away entirely; clear their
# formatting first so no stale style reference to them survives, then wipe
# every cell and start the layout fresh.
$ws.Columns.Item(5).ClearFormats()
$ws.Columns.Item(4).ClearFormats()
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "NOMBRE"
$ws.Range("B1").Value = "CAPACITACION "
$ws.Range("C1").Value = "URL"

# Four identical data rows
$nombre = "ALBA AIDEE AGUILAR ABREGO"
$capacitacion = "AUDITOR INTERNO DE SISTEMAS DE GESTIÓN DE LA CALIDAD SEGÚN LA NORMA ISO 9001:2015"
$url = "https://servimeters-my.sharepoint.com/:x:/r/personal/duvan_sanabria_servimeters_com/_layouts/15/Doc.aspx?sourcedoc=%7B9D631C9C-C83D-4260-8520-0BA5934135E4%7D&file=Prog%20en%20Sitio-RICARDO.xlsx&action=default&mobileredirect=true&DefaultItemOpen=1"

for ($r = 2; $r -le 5; $r++) {
    $ws.Range("A$r").Value = $nombre
    $ws.Range("B$r").Value = $capacitacion
    $ws.Range("C$r").Value = $url
}

# Column C is text-formatted throughout; rows 2-5 additionally pick up the
# "Hipervínculo" look (blue/underline) even though the live hyperlink is gone.
$ws.Range("C1:C5").NumberFormat = "@"
$ws.Range("C2:C5").Style = "Hipervínculo"
$ws.Range("C2:C5").NumberFormat = "@"

# New column widths for the three surviving columns.
$ws.Columns.Item(1).ColumnWidth = 31.85546875
$ws.Columns.Item(2).ColumnWidth = 90.140625
$ws.Columns.Item(3).ColumnWidth = 38

# Selection moves to the new data block.
$ws.Range("A3:C5").Select()
